$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list update: refresh Price (D) and Volume(1h) (E) columns.
# D-column values are stored as literal text (e.g. "10.50", "57.050.33"),
# so force NumberFormat to Text before writing, then restore the default
# "Normal" style so no stray formatting is left behind on the cell.

# Row 2
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '56.752.48'
$c.Style = "Normal"
$ws.Range("E2").Value = '  +3.40%  '

# Row 3
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '2.329.45'
$c.Style = "Normal"
$ws.Range("E3").Value = '  +1.53%  '

# Row 4
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.01'
$c.Style = "Normal"
$ws.Range("E4").Value = '  +0.48%  '

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '520.43'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +2.52%  '

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '134.72'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +3.80%  '

# Row 7
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.997'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +0.24%  '

# Row 8
$ws.Range("E8").Value = '  +1.57%  '

# Row 9
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '2.356.12'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +1.69%  '

# Row 10
$ws.Range("E10").Value = '  +6.73%  '

# Row 11
$ws.Range("E11").Value = '  -0.79%  '

# Row 12
$ws.Range("E12").Value = '  +3.58%  '

# Row 13
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '0.342'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +0.59%  '

# Row 14
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '23.77'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -0.52%  '

# Row 15
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '2.750.47'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +1.65%  '

# Row 16
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '56.861.54'
$c.Style = "Normal"
$ws.Range("E16").Value = '  +3.63%  '

# Row 17
$ws.Range("E17").Value = '  +2.39%  '

# Row 18
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '2.353.39'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +3.52%  '

# Row 19
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '10.50'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -1.86%  '

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '324.54'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +4.73%  '

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '6.54'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -2.38%  '

# Row 23
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '0.998'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +0.03%  '

# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '60.95'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +0.59%  '

# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '0.163'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +7.52%  '

# Row 26
$ws.Range("E26").Value = '  +0.42%  '

# Row 27
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '7.87'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +4.59%  '

# Row 28
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '1.27'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +10.41%  '

# Row 29
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '170.84'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -1.08%  '

# Row 30
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '0.0₃0748'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +5.48%  '

# Row 31
$ws.Range("E31").Value = '  +3.72%  '

# Row 32
$ws.Range("E32").Value = '  +0.64%  '

# Row 33
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '18.32'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +1.16%  '

# Row 34
$ws.Range("E34").Value = '  +0.04%  '

# Row 35
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.992'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -0.02%  '

# Row 36
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '1.25'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +1.86%  '

# Row 37
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.923'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +0.84%  '

# Row 38
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '4.03'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +3.71%  '

# Row 39
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '1.56'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +8.09%  '

# Row 40
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '37.91'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +3.12%  '

# Row 41
$ws.Range("E41").Value = '  +0.24%  '

# Row 42
$ws.Range("E42").Value = '  +4.42%  '

# Row 43
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '137.50'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +1.35%  '

# Row 44
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '279.69'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +8.56%  '

# Row 45
$ws.Range("E45").Value = '  +0.43%  '

# Row 46
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.0936'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +2.58%  '

# Row 47
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.0505'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +0.22%  '

# Row 48
$ws.Range("E48").Value = '  +1.92%  '

# Row 49
$ws.Range("E49").Value = '  +4.52%  '

# Row 50
$ws.Range("E50").Value = '  +0.42%  '

# Row 51
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '17.59'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +6.44%  '
